$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing two columns (old A,B) right into B,C so a new column A
# ("Key") can be introduced in front of them; the old C column becomes D.
$ws.Columns("A").Insert()

# Make room for two new rows of data (inserted after the current row 4,
# i.e. becoming the new rows 5 and 6) describing additional categories.
$ws.Rows("5:6").Insert()

# ---- Header row ----
$ws.Cells.Item(1,1).Value = "Key"
$ws.Cells.Item(1,2).Value = "Main category"
$ws.Cells.Item(1,3).Value = "Sub-category"
$ws.Cells.Item(1,4).Value = "Example"

# ---- Row 2 ----
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = "Not applicable (protocol for a review, commentary, etc)"
$ws.Cells.Item(2,4).Value = '"Data sharing not applicable to this article as no datasets were generated or analysed during the current study."[@ehrlich2019]'

# ---- Row 3 ----
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = "Data not available"
$ws.Cells.Item(3,3).Value = "Data not made available"
$ws.Cells.Item(3,4).Value = '"Not available for public"[@septiandri2019]'

# ---- Row 4 ----
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = "Data not available"
$ws.Cells.Item(4,3).Value = "Data available on request to authors"
$ws.Cells.Item(4,4).Value = '"Data can be available upon reasonable request to the corresponding author."[@solis2019]'

# ---- Row 5 (new) ----
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = "Data not available"
$ws.Cells.Item(5,3).Value = "Data will be available in the future"
$ws.Cells.Item(5,4).Value = '"The protocol and full dataset will be available at Open Science Framework upon peer review publication (https://osf.io/rvbuy/)."[@ebbeling2019]'

# ---- Row 6 (new) ----
$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = "Data not available"
$ws.Cells.Item(6,3).Value = "Data vailable from central repository, but insufficient detail published to find"
$ws.Cells.Item(6,4).Value = '"Data were obtained from the international MSBase cohort study. Information regarding data availability can be obtained at https://www.msbase.org/."[@malpas2019]'

# ---- Row 7 ----
$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = "Data available"
$ws.Cells.Item(7,3).Value = "Data available in the manuscript/supplementary files"
$ws.Cells.Item(7,4).Value = '"All data related to this study are present in the paper or the Supplementary Materials. . ."[@thompson2019]'

# ---- Row 8 ----
$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = "Data available"
$ws.Cells.Item(8,3).Value = "Data available in online repository e.g. GitHub, Zenodo"
$ws.Cells.Item(8,4).Value = '"Extracted data used in this meta-analysis and analysis code are available at www.doi.org/10.5281/zenodo.3149365."[@moriarty2019a]'

# ---- Row 9 ----
$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = "Data available"
$ws.Cells.Item(9,3).Value = "Data available from central repository (requires sufficient details to identify e.g. extract or accession ID)"
$ws.Cells.Item(9,4).Value = '"This research has been conducted using the UK Biobank Resource under application number 24494. All bona fide researchers can apply to use the UK Biobank resource for health related research that is in the public interest."[@knuppel2019]'

# Empty, word-wrapped placeholder cell further down the sheet.
$ws.Range("C14").WrapText = $true

# Match the saved selection / active cell and printer page setup.
[void]$ws.Range("C9").Select()
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
